# Padel workbook update
# - resultados: add new match result row (Los luises vs Ceci-Angel, "Mediocre alto" group)
# - clasificacion_auto / clasificacion: update standings for "Los luises" (winner) and
#   "Ceci-Ángel" (loser) to reflect the new match
# - historial_partidos: append the two match-history ledger rows for the new match

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) resultados -> new row with the match result
# ---------------------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("resultados")
$wsRes.Cells.Item(23, 1).Value = "Mediocre alto"
$wsRes.Cells.Item(23, 2).Value = "1ª vuelta"
$wsRes.Cells.Item(23, 3).Value = "Los luises"
$wsRes.Cells.Item(23, 4).Value = "Ceci-Ángel"
$wsRes.Cells.Item(23, 5).Value = "6-1,6-3"
$wsRes.Cells.Item(23, 6).Value = "1-6,3-6"

$wsRes.Activate()
$wsRes.Range("F23").Select()

# ---------------------------------------------------------------------------
# 2) clasificacion_auto -> update accumulated stats
#    row 2 = Los luises (mediocre alto), row 6 = Ceci-Angel (mediocre alto)
# ---------------------------------------------------------------------------
$wsAuto = $wb.Worksheets.Item("clasificacion_auto")

# Los luises: PUNTOS, PJ, PG, SG, JG, JP
$wsAuto.Cells.Item(2, 4).Value = 10
$wsAuto.Cells.Item(2, 5).Value = 4
$wsAuto.Cells.Item(2, 6).Value = 3
$wsAuto.Cells.Item(2, 9).Value = 7
$wsAuto.Cells.Item(2, 11).Value = 53
$wsAuto.Cells.Item(2, 12).Value = 35

# Ceci-Angel: PJ, PP, SP, JG, JP
$wsAuto.Cells.Item(6, 5).Value = 3
$wsAuto.Cells.Item(6, 8).Value = 3
$wsAuto.Cells.Item(6, 10).Value = 6
$wsAuto.Cells.Item(6, 11).Value = 12
$wsAuto.Cells.Item(6, 12).Value = 37

# ---------------------------------------------------------------------------
# 3) clasificacion -> update accumulated stats (no JG/JP columns here)
#    row 2 = Los luises (Mediocre alto), row 6 = Ceci-Angel (Mediocre alto)
# ---------------------------------------------------------------------------
$wsClas = $wb.Worksheets.Item("clasificacion")

# Los luises: PUNTOS, P.JUGADOS, P.GANADOS, SET GANADOS
$wsClas.Cells.Item(2, 4).Value = 10
$wsClas.Cells.Item(2, 5).Value = 4
$wsClas.Cells.Item(2, 6).Value = 3
$wsClas.Cells.Item(2, 9).Value = 7

# Ceci-Angel: P.JUGADOS, P.PERDIDOS, SET PERDIDOS
$wsClas.Cells.Item(6, 5).Value = 3
$wsClas.Cells.Item(6, 8).Value = 3
$wsClas.Cells.Item(6, 10).Value = 6

# ---------------------------------------------------------------------------
# 4) historial_partidos -> append ledger rows 44 (winner) & 45 (loser)
# ---------------------------------------------------------------------------
$wsHist = $wb.Worksheets.Item("historial_partidos")

# copy the formatting of the last existing row (43) down onto the two new rows
# so the FECHA column keeps its date number format (style index) intact
$wsHist.Range("A43:M43").Copy()
$wsHist.Range("A44:M45").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 44: Los luises win
$wsHist.Cells.Item(44, 1).Value = 45974
$wsHist.Cells.Item(44, 2).Value = "mediocre alto"
$wsHist.Cells.Item(44, 3).Value = "1ª vuelta"
$wsHist.Cells.Item(44, 4).Value = "Los luises"
$wsHist.Cells.Item(44, 5).Value = "Gana"
$wsHist.Cells.Item(44, 6).Value = 2
$wsHist.Cells.Item(44, 7).Value = 0
$wsHist.Cells.Item(44, 8).Value = 3
$wsHist.Cells.Item(44, 9).Value = 4
$wsHist.Cells.Item(44, 10).Value = 10
$wsHist.Cells.Item(44, 11).Value = 3
$wsHist.Cells.Item(44, 12).Value = 1
$wsHist.Cells.Item(44, 13).Value = 0

# Row 45: Ceci-Angel loss
$wsHist.Cells.Item(45, 1).Value = 45974
$wsHist.Cells.Item(45, 2).Value = "mediocre alto"
$wsHist.Cells.Item(45, 3).Value = "1ª vuelta"
$wsHist.Cells.Item(45, 4).Value = "Ceci-Ángel"
$wsHist.Cells.Item(45, 5).Value = "Pierde"
$wsHist.Cells.Item(45, 6).Value = 0
$wsHist.Cells.Item(45, 7).Value = 2
$wsHist.Cells.Item(45, 8).Value = 0
$wsHist.Cells.Item(45, 9).Value = 3
$wsHist.Cells.Item(45, 10).Value = 0
$wsHist.Cells.Item(45, 11).Value = 0
$wsHist.Cells.Item(45, 12).Value = 0
$wsHist.Cells.Item(45, 13).Value = 3
